# Day5-second commit: automation adds patient details for the appointment.
# Fills in First Name / Last Name / email columns next to the existing
# "number" column, turns the email into a mailto: hyperlink (Excel's
# built-in "Hyperlink" cell style), widens the new columns to fit their
# content, and leaves the selection on B3 (as if the operator had just
# tabbed/entered past the newly-typed row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (B1:D1)
$ws.Range("B1").Value = "First Name"
$ws.Range("C1").Value = "Last Name"
$ws.Range("D1").Value = "email"

# New patient row (B2:D2) alongside the existing phone number in A2
$ws.Range("B2").Value = "abi"
$ws.Range("C2").Value = "s"
$ws.Range("D2").Value = "abi@gmail.com"

# Turn the email address into a clickable mailto hyperlink (applies the
# built-in "Hyperlink" style - underline + theme color - to D2)
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:abi@gmail.com") | Out-Null

# Widen the new columns so the patient details are fully visible
$ws.Columns.Item(2).ColumnWidth = 12.666666666666666
$ws.Columns.Item(3).ColumnWidth = 13.333333333333332
$ws.Columns.Item(4).ColumnWidth = 14.333333333333332
$ws.Columns.Item(5).ColumnWidth = 16.333333333333336

# Leave the selection where the operator's typing ended up
$ws.Range("B3").Select() | Out-Null
